$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values while preserving their original text cell type
# (quote-prefix forces text interpretation; Style reset avoids leaving a
# stray "quote prefix" style applied to the cell).

$ws.Range("D2").Value = "'246.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.27%"
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").Value = "'11"
$ws.Range("G2").Style = "Normal"
$ws.Range("D3").Value = "'30.13"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-0.16%"
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").Value = "'11"
$ws.Range("G3").Style = "Normal"
$ws.Range("D4").Value = "'5.159"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.19%"
$ws.Range("E4").Style = "Normal"
$ws.Range("G4").Value = "'11"
$ws.Range("G4").Style = "Normal"
$ws.Range("D5").Value = "'0.05759"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.48%"
$ws.Range("E5").Style = "Normal"
$ws.Range("G5").Value = "'11"
$ws.Range("G5").Style = "Normal"
$ws.Range("D6").Value = "'6.668"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.01%"
$ws.Range("E6").Style = "Normal"
$ws.Range("G6").Value = "'11"
$ws.Range("G6").Style = "Normal"
$ws.Range("D7").Value = "'3.265"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'6.41%"
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").Value = "'11"
$ws.Range("G7").Style = "Normal"
$ws.Range("D8").Value = "'0.8495"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-1.07%"
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").Value = "'11"
$ws.Range("G8").Style = "Normal"
$ws.Range("D9").Value = "'0.8571"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-2.57%"
$ws.Range("E9").Style = "Normal"
$ws.Range("G9").Value = "'11"
$ws.Range("G9").Style = "Normal"
$ws.Range("D10").Value = "'0.1389"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.57%"
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").Value = "'11"
$ws.Range("G10").Style = "Normal"
$ws.Range("D11").Value = "'0.07088"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'0.15%"
$ws.Range("E11").Style = "Normal"
$ws.Range("G11").Value = "'11"
$ws.Range("G11").Style = "Normal"
$ws.Range("D12").Value = "'0.03264"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'13.87%"
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").Value = "'11"
$ws.Range("G12").Style = "Normal"
$ws.Range("D13").Value = "'0.09370"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.21%"
$ws.Range("E13").Style = "Normal"
$ws.Range("G13").Value = "'11"
$ws.Range("G13").Style = "Normal"
$ws.Range("D14").Value = "'0.001529"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'1.30%"
$ws.Range("E14").Style = "Normal"
$ws.Range("G14").Value = "'11"
$ws.Range("G14").Style = "Normal"
$ws.Range("D15").Value = "'0.0005980"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-94.19%"
$ws.Range("E15").Style = "Normal"
$ws.Range("G15").Value = "'11"
$ws.Range("G15").Style = "Normal"
$ws.Range("D16").Value = "'0.005906"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.75%"
$ws.Range("E16").Style = "Normal"
$ws.Range("G16").Value = "'11"
$ws.Range("G16").Style = "Normal"
$ws.Range("D17").Value = "'3.528"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.95%"
$ws.Range("E17").Style = "Normal"
$ws.Range("G17").Value = "'11"
$ws.Range("G17").Style = "Normal"
$ws.Range("D18").Value = "'2.218"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'1.79%"
$ws.Range("E18").Style = "Normal"
$ws.Range("G18").Value = "'11"
$ws.Range("G18").Style = "Normal"
$ws.Range("D19").Value = "'0.3123"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.94%"
$ws.Range("E19").Style = "Normal"
$ws.Range("G19").Value = "'11"
$ws.Range("G19").Style = "Normal"
$ws.Range("D20").Value = "'0.03415"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'4.57%"
$ws.Range("E20").Style = "Normal"
$ws.Range("G20").Value = "'11"
$ws.Range("G20").Style = "Normal"
$ws.Range("D21").Value = "'0.1316"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.17%"
$ws.Range("E21").Style = "Normal"
$ws.Range("G21").Value = "'11"
$ws.Range("G21").Style = "Normal"
$ws.Range("D22").Value = "'3.499"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-1.05%"
$ws.Range("E22").Style = "Normal"
$ws.Range("G22").Value = "'11"
$ws.Range("G22").Style = "Normal"
$ws.Range("B23").Value = "'ZBToken"
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").Value = "'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").Value = "'0.1410"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'2.21%"
$ws.Range("E23").Style = "Normal"
$ws.Range("G23").Value = "'11"
$ws.Range("G23").Style = "Normal"
$ws.Range("B24").Value = "'CoinExToken"
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").Value = "'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").Value = "'0.04121"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.74%"
$ws.Range("E24").Style = "Normal"
$ws.Range("G24").Value = "'11"
$ws.Range("G24").Style = "Normal"
$ws.Range("E25").Value = "'1.07%"
$ws.Range("E25").Style = "Normal"
$ws.Range("G25").Value = "'11"
$ws.Range("G25").Style = "Normal"
$ws.Range("D26").Value = "'0.004158"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-7.61%"
$ws.Range("E26").Style = "Normal"
$ws.Range("G26").Value = "'11"
$ws.Range("G26").Style = "Normal"
$ws.Range("E27").Value = "'-0.83%"
$ws.Range("E27").Style = "Normal"
$ws.Range("G27").Value = "'11"
$ws.Range("G27").Style = "Normal"
$ws.Range("D28").Value = "'0.0001449"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'4.77%"
$ws.Range("E28").Style = "Normal"
$ws.Range("G28").Value = "'11"
$ws.Range("G28").Style = "Normal"
$ws.Range("G29").Value = "'11"
$ws.Range("G29").Style = "Normal"
$ws.Range("G30").Value = "'11"
$ws.Range("G30").Style = "Normal"
$ws.Range("G31").Value = "'11"
$ws.Range("G31").Style = "Normal"
$ws.Range("G32").Value = "'11"
$ws.Range("G32").Style = "Normal"
$ws.Range("G33").Value = "'11"
$ws.Range("G33").Style = "Normal"
$ws.Range("G34").Value = "'11"
$ws.Range("G34").Style = "Normal"
$ws.Range("G35").Value = "'11"
$ws.Range("G35").Style = "Normal"
$ws.Range("G36").Value = "'11"
$ws.Range("G36").Style = "Normal"
$ws.Range("G37").Value = "'11"
$ws.Range("G37").Style = "Normal"
$ws.Range("G38").Value = "'11"
$ws.Range("G38").Style = "Normal"
$ws.Range("G39").Value = "'11"
$ws.Range("G39").Style = "Normal"
$ws.Range("D40").Value = "'0.03751"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.97%"
$ws.Range("E40").Style = "Normal"
$ws.Range("G40").Value = "'11"
$ws.Range("G40").Style = "Normal"
$ws.Range("D41").Value = "'0.1072"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.01%"
$ws.Range("E41").Style = "Normal"
$ws.Range("G41").Value = "'11"
$ws.Range("G41").Style = "Normal"
$ws.Range("D42").Value = "'0.002460"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-5.38%"
$ws.Range("E42").Style = "Normal"
$ws.Range("G42").Value = "'11"
$ws.Range("G42").Style = "Normal"
$ws.Range("D43").Value = "'0.003544"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-38.14%"
$ws.Range("E43").Style = "Normal"
$ws.Range("G43").Value = "'11"
$ws.Range("G43").Style = "Normal"
$ws.Range("D44").Value = "'0.01037"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'3.24%"
$ws.Range("E44").Style = "Normal"
$ws.Range("G44").Value = "'11"
$ws.Range("G44").Style = "Normal"
$ws.Range("E45").Value = "'7.64%"
$ws.Range("E45").Style = "Normal"
$ws.Range("G45").Value = "'11"
$ws.Range("G45").Style = "Normal"
$ws.Range("E46").Value = "'0.02%"
$ws.Range("E46").Style = "Normal"
$ws.Range("G46").Value = "'11"
$ws.Range("G46").Style = "Normal"
$ws.Range("D47").Value = "'0.07099"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-20.21%"
$ws.Range("E47").Style = "Normal"
$ws.Range("G47").Value = "'11"
$ws.Range("G47").Style = "Normal"
$ws.Range("E48").Value = "'-10.85%"
$ws.Range("E48").Style = "Normal"
$ws.Range("G48").Value = "'11"
$ws.Range("G48").Style = "Normal"
$ws.Range("G49").Value = "'11"
$ws.Range("G49").Style = "Normal"
$ws.Range("E50").Value = "'0.02%"
$ws.Range("E50").Style = "Normal"
$ws.Range("G50").Value = "'11"
$ws.Range("G50").Style = "Normal"
$ws.Range("G51").Value = "'11"
$ws.Range("G51").Style = "Normal"
